# Updates the cryptos price/volume table to the latest snapshot
# (GitHub Actions scheduled refresh).
#
# Column D ("Price") values look numeric (e.g. "51.706.52", "1.00",
# "0.300") but are really formatted text straight from the source feed
# -- grouped with dots, and sometimes carrying trailing zeros that a
# real number would drop. Each is prefixed with a leading apostrophe,
# exactly like typing '51.70 directly into Excel, so it is stored
# verbatim as text instead of being re-parsed/rounded as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'51.706.52"
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = "'3.036.07"
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'380.94"
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').Value = "'102.92"
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('E7').Value = '  +1.12%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').Value = "'36.82"
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = "'0.0866"
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').Value = "'3.512.48"
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').Value = "'7.75"
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('D16').Value = "'3.026.50"
$ws.Range('E16').Value = '  +2.27%  '
$ws.Range('D17').Value = "'0.978"
$ws.Range('E17').Value = '  -3.77%  '
$ws.Range('D18').Value = "'10.49"
$ws.Range('E18').Value = '  -15.13%  '
$ws.Range('D19').Value = "'51.744.58"
$ws.Range('E19').Value = '  +1.40%  '
$ws.Range('D20').Value = "'3.10"
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('D23').Value = "'70.08"
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('D24').Value = "'267.76"
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('E25').Value = '  -6.34%  '
$ws.Range('D26').Value = "'8.28"
$ws.Range('E26').Value = '  +2.80%  '
$ws.Range('D27').Value = "'7.67"
$ws.Range('E27').Value = '  +9.98%  '
$ws.Range('E28').Value = '  +4.81%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = "'26.26"
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('D32').Value = "'10.29"
$ws.Range('E32').Value = '  -1.03%  '
$ws.Range('D33').Value = "'2.08"
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('D35').Value = "'33.86"
$ws.Range('E35').Value = '  -0.66%  '
$ws.Range('D36').Value = "'0.0447"
$ws.Range('E36').Value = '  +2.81%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').Value = "'3.32"
$ws.Range('E38').Value = '  +4.08%  '
$ws.Range('D39').Value = "'0.300"
$ws.Range('E39').Value = '  +17.11%  '
$ws.Range('D40').Value = "'17.04"
$ws.Range('E40').Value = '  +2.36%  '
$ws.Range('E41').Value = '  +2.32%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = "'127.75"
$ws.Range('E42').Value = '  +7.80%  '
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = "'2.55"
$ws.Range('E44').Value = '  +2.45%  '
$ws.Range('D45').Value = "'3.78"
$ws.Range('E45').Value = '  +6.06%  '
$ws.Range('D46').Value = "'21.72"
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('D47').Value = "'2.10"
$ws.Range('E47').Value = '  +3.60%  '
$ws.Range('D48').Value = "'2.40"
$ws.Range('E48').Value = '  +3.66%  '
$ws.Range('D49').Value = "'2.031.43"
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = "'3.336.89"
$ws.Range('E50').Value = '  +2.66%  '
$ws.Range('E51').Value = '  +1.27%  '
